# Apply the "Country of Birth added to patient profile" edit.
#
# This workbook documents a FHIR StructureDefinition ("historical-data")
# as two sheets: "Metadata" and "Elements".
#
# Changes:
#  1. Metadata!B8 - the "Date" metadata value is bumped to the new commit
#     timestamp.
#  2. Elements sheet, row for "Extension.extension" (row 3): the slice
#     is turned into a closed-off slot (max cardinality 0) with a plain
#     "Extension" / "An Extension" short label + definition, and the
#     slicing/comment/requirements text is cleared out.
#  3. Elements sheet, row for "Extension.value[x]" (row 5): the allowed
#     type list is collapsed down to just "string".

$wb = $excel.ActiveWorkbook

# ---- 1. Metadata sheet: update Date value ----
$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("B8").Value = "2023-07-13T10:52:26+02:00"

# ---- 2 & 3. Elements sheet ----
$elements = $wb.Worksheets.Item("Elements")

# Row 3 = Extension.extension
# D3 and AK3 become empty strings (still "string typed" blank cells, like
# the other already-blank cells in the row); copy from one of those
# existing blank cells so the stored cell keeps its shared-string type.
$elements.Range("H3").Copy($elements.Range("D3"))
$elements.Range("H3").Copy($elements.Range("AK3"))
$elements.Range("G3").Value = "0"
$elements.Range("L3").Value = "Extension"
$elements.Range("M3").Value = "An Extension"
# N3 is cleared out entirely (no value at all).
$elements.Range("N3").ClearContents()

# Row 5 = Extension.value[x]
$elements.Range("K5").Value = "string`n"
